$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header names / account number
$ws.Range("C2").Value = "Hartmut"

# B3 holds a 16-digit card number that must stay TEXT (it would otherwise
# be coerced to a Number, rendered in scientific notation / lose digits
# under Excel's 15-significant-digit rule). Force text entry by switching
# to a text number format, assigning the value, then pasting the original
# cell's format back on top so the style index is unchanged (still s="8")
# while the stored cell keeps its new string type/value.
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("C3").Copy()
$ws.Range("B3").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("C3").Value = "Mohaupt"

# Opening balance line
$ws.Range("D5").Value = "KONTOSTAND AM 01.02.2024"

# Row 6
$ws.Range("B6").Value = "03.02."
$ws.Range("C6").Value = "04.02."
$ws.Range("D6").Value = "PAYPAL VVPWAU"
$ws.Range("E6").Value = "44,78-"

# Row 7
$ws.Range("B7").Value = "05.02."
$ws.Range("C7").Value = "06.02."
$ws.Range("D7").Value = "KARTENZ./05.02 ALDI SUED RO"
$ws.Range("E7").Value = "128,97-"

# Row 8
$ws.Range("B8").Value = "07.02."
$ws.Range("C8").Value = "08.02."
$ws.Range("D8").Value = "EBAY MKTPLC EU LMZWXN"
$ws.Range("E8").Value = "110,61-"

# Row 9
$ws.Range("B9").Value = "10.02."
$ws.Range("C9").Value = "11.02."
$ws.Range("D9").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 11805213"
$ws.Range("E9").Value = "87,59-"

# Row 10 - was empty, now filled in. Copy style from row 9 so the new
# entries look like the other transaction rows (E column needs the
# right-aligned "amount" style used by E6:E9, not the blank E10/E11 style).
$ws.Range("B9:E9").Copy()
$ws.Range("B10:E10").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B10").Value = "13.02."
$ws.Range("C10").Value = "14.02."
$ws.Range("D10").Value = "KARTENZAHLUNG ARAL TANKSTELLE"
$ws.Range("E10").Value = "89,50-"

# Row 11
$ws.Range("B9:E9").Copy()
$ws.Range("B11:E11").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B11").Value = "15.02."
$ws.Range("C11").Value = "16.02."
$ws.Range("D11").Value = "PAYPAL SHQWRP"
$ws.Range("E11").Value = "60,08-"

$excel.CutCopyMode = $false

# Closing balance line
$ws.Range("D12").Value = "KONTOSTAND AM 19.02.2024"
$ws.Range("E12").Value = "521,53-"

# Next billing date note
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 28.02.2024"
